# Update the "datetimeFigureOut" date field shown in the footer area of
# every slide layout and the slide master (PowerPoint recomputes/caches
# this field's displayed text; here we bump it from 12/20/2025 to
# 12/21/2025, matching the authored change).
#
# Note: helper functions in this COM shim must take positional
# parameters (named parameters like `-shapes` lose the live COM
# reference), so Set-DateFieldText below is called positionally.

$p = $ppt.ActivePresentation

function Set-DateFieldText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

$newDate = "12/21/2025"

$master = $p.SlideMaster
Set-DateFieldText $master.Shapes $newDate

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Set-DateFieldText $layout.Shapes $newDate
}

# Slide 20: update the Tableau dashboard link text (new workbook/view name
# and a trimmed query string with the "publish=yes" param removed). Find
# the shape by its current text (rather than a hard-coded index) so the
# edit is resilient to any shape re-ordering.
$newTableauUrl = "https://public.tableau.com/views/Air_Craft_final_project/Dashboard1?:language=en-US&:sid=&:redirect=auth&:display_count=n&:origin=viz_share_link"
$s20 = $p.Slides.Item(20)
for ($i = 1; $i -le $s20.Shapes.Count; $i++) {
    $shp = $s20.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -like "*public.tableau.com*") {
        $shp.TextFrame.TextRange.Text = $newTableauUrl
    }
}
